# Auto-generated edit script for LOM3223.docx
$d = $word.ActiveDocument
$nl = [char]11

# Para 6 run 0
$rng = $d.Paragraphs.Item(6).Range
$ok = $rng.Find.Execute('Propiciar ao aluno os conhecimentos básicos de materiais magnéticos e supercondutores visando sua aplicação em dispositivos.', $true, $false, $false, $false, $false, $true, 1, $false, 'Conceitos fundamentais de propriedades magnéticas da matéria. Magnetismo de elétrons. Ferromagnetismo. Materiais magnéticos e Aplicações: moles e duros. Interação de troca em óxidos e metais. Magnetismo - Fenomenologia Clássica: diamagnetismo e paramagnetismo. Magnetismo - Fenomenologia Quântica: ferromagnetismo. Anisotropia Magnética e Interação Spin-Órbita. Magnetostricção e materiais magnetostrictivos -Introdução e aplicações. Conceitos básicos de supercondutividade. Supercondutividade - Origem Quântica. Super-onda – Consequências. Interferencia quântica – SQUID. Materiais Supercondutores e Aplicações.', 2)
if (-not $ok) { Write-Output "FAILED: para 6 run 0" } else { Write-Output "OK: para 6 run 0" }

# Para 7 run 0
$rng = $d.Paragraphs.Item(7).Range
$ok = $rng.Find.Execute('Provide the student with the basic knowledge of magnetic and superconducting materials aiming their application in devices.', $true, $false, $false, $false, $false, $true, 1, $false, 'Fundamental concepts of magnetic properties of matter. Electron magnetism. Ferromagnetism. Magnetic Materials and Applications: soft and hard. Exchange interaction in oxides and metals. Magnetism - Classical Phenomenology: diamagnetism and paramagnetism. Magnetism - Quantum Phenomenology: ferromagnetism. Magnetic Anisotropy and Spin-Orbit Interaction. Magnetostriction and magnetostrictive materials -Introduction and applications. Basic concepts of superconductivity. Superconductivity - Quantum Origin. Superwave – Consequences. Quantum Interference – SQUID. Superconducting Materials and Applications', 2)
if (-not $ok) { Write-Output "FAILED: para 7 run 0" } else { Write-Output "OK: para 7 run 0" }

# Para 9 run 0
$rng = $d.Paragraphs.Item(9).Range
$ok = $rng.Find.Execute('5840730 - Antonio Jefferson da Silva Machado' + $nl + '', $true, $false, $false, $false, $false, $true, 1, $false, 'Propiciar ao aluno os conhecimentos básicos de materiais magnéticos e supercondutores visando sua aplicação em dispositivos.' + $nl + '', 2)
if (-not $ok) { Write-Output "FAILED: para 9 run 0" } else { Write-Output "OK: para 9 run 0" }

# Para 9 run 1
$rng = $d.Paragraphs.Item(9).Range
$ok = $rng.Find.Execute('5840726 - Cristina Bormio Nunes', $true, $false, $false, $false, $false, $true, 1, $false, 'Conceitos fundamentais de propriedades magnéticas da matéria. Magnetismo de elétrons. Ferromagnetismo. Materiais magnéticos e Aplicações: moles e duros. Interação de troca em óxidos e metais. Magnetismo - Fenomenologia Clássica: diamagnetismo e paramagnetismo. Magnetismo - Fenomenologia Quântica: ferromagnetismo. Anisotropia Magnética e Interação Spin-Órbita. Magnetostricção e materiais magnetostrictivos -Introdução e aplicações. Conceitos básicos de supercondutividade. Supercondutividade - Origem Quântica. Super-onda – Consequências. Interferencia quântica – SQUID. Materiais Supercondutores e Aplicações.', 2)
if (-not $ok) { Write-Output "FAILED: para 9 run 1" } else { Write-Output "OK: para 9 run 1" }

# Para 11 run 0
$rng = $d.Paragraphs.Item(11).Range
$ok = $rng.Find.Execute('Conceitos fundamentais de propriedades magnéticas da matéria. Magnetismo de elétrons. Ferromagnetismo. Materiais magnéticos e Aplicações: moles e duros. Interação de troca em óxidos e metais. Magnetismo - Fenomenologia Clássica: diamagnetismo e paramagnetismo. Magnetismo - Fenomenologia Quântica: ferromagnetismo. Anisotropia Magnética e Interação Spin-Órbita. Magnetostricção e materiais magnetostrictivos -Introdução e aplicações. Conceitos básicos de supercondutividade. Supercondutividade - Origem Quântica. Super-onda – Consequências. Interferencia quântica – SQUID. Materiais Supercondutores e Aplicações.', $true, $false, $false, $false, $false, $true, 1, $false, 'Aulas expositivas, seminários e exercícios comentados.', 2)
if (-not $ok) { Write-Output "FAILED: para 11 run 0" } else { Write-Output "OK: para 11 run 0" }

# Para 12 run 0
$rng = $d.Paragraphs.Item(12).Range
$ok = $rng.Find.Execute('Fundamental concepts of magnetic properties of matter. Electron magnetism. Ferromagnetism. Magnetic Materials and Applications: soft and hard. Exchange interaction in oxides and metals. Magnetism - Classical Phenomenology: diamagnetism and paramagnetism. Magnetism - Quantum Phenomenology: ferromagnetism. Magnetic Anisotropy and Spin-Orbit Interaction. Magnetostriction and magnetostrictive materials -Introduction and applications. Basic concepts of superconductivity. Superconductivity - Quantum Origin. Superwave – Consequences. Quantum Interference – SQUID. Superconducting Materials and Applications', $true, $false, $false, $false, $false, $true, 1, $false, 'Provide the student with the basic knowledge of magnetic and superconducting materials aiming their application in devices.', 2)
if (-not $ok) { Write-Output "FAILED: para 12 run 0" } else { Write-Output "OK: para 12 run 0" }

# Para 14 run 0
$rng = $d.Paragraphs.Item(14).Range
$ok = $rng.Find.Execute('Conceitos fundamentais de propriedades magnéticas da matéria. Magnetismo de elétrons. Ferromagnetismo. Materiais magnéticos e Aplicações: moles e duros. Interação de troca em óxidos e metais. Magnetismo - Fenomenologia Clássica: diamagnetismo e paramagnetismo. Magnetismo - Fenomenologia Quântica: ferromagnetismo. Anisotropia Magnética e Interação Spin-Órbita. Magnetostricção e materiais magnetostrictivos -Introdução e aplicações. Conceitos básicos de supercondutividade. Supercondutividade - Origem Quântica. Super-onda – Consequências. Interferencia quântica – SQUID. Materiais Supercondutores e Aplicações.', $true, $false, $false, $false, $false, $true, 1, $false, 'A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável.', 2)
if (-not $ok) { Write-Output "FAILED: para 14 run 0" } else { Write-Output "OK: para 14 run 0" }

# Para 17 run 1
$rng = $d.Paragraphs.Item(17).Range
$ok = $rng.Find.Execute('Aulas expositivas, seminários e exercícios comentados.' + $nl + '', $true, $false, $false, $false, $false, $true, 1, $false, 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação' + $nl + '', 2)
if (-not $ok) { Write-Output "FAILED: para 17 run 1" } else { Write-Output "OK: para 17 run 1" }

# Para 17 run 3
$rng = $d.Paragraphs.Item(17).Range
$ok = $rng.Find.Execute('A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável.' + $nl + '', $true, $false, $false, $false, $false, $true, 1, $false, 'JILES, D. C. Introduction to Magnetism and Magnetic Materials, CRC Press, 1998.' + $nl + 'COEY, J. M. D. Magnetism and Magnetic Materials, Cambridge University Press, 2010.' + $nl + 'BUSCHOW, K. H. J.; DE BOER, F. R. Physics of Magnetism and Magnetic Materials, Springer, 2003.' + $nl + 'CULLITY, B. D.; GRAHAM, C. D. Introduction to Magnetic Materials, Wiley-IEEE Press, 2008.' + $nl + 'POOLE, C. P. et al., Superconductivity, Academic Press, 2007.' + $nl + 'SHEAHEN, T. P. Introduction to High-Temperature Superconductivity, Kluwer Academic, 2002.' + $nl + 'LEE, P. J. Engineering Superconductivity, Wiley-IEEE Press, 2001.' + $nl + '', 2)
if (-not $ok) { Write-Output "FAILED: para 17 run 3" } else { Write-Output "OK: para 17 run 3" }

# Para 17 run 5
$rng = $d.Paragraphs.Item(17).Range
$ok = $rng.Find.Execute('Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação', $true, $false, $false, $false, $false, $true, 1, $false, '5840730 - Antonio Jefferson da Silva Machado', 2)
if (-not $ok) { Write-Output "FAILED: para 17 run 5" } else { Write-Output "OK: para 17 run 5" }

# Para 19 run 0
$rng = $d.Paragraphs.Item(19).Range
$ok = $rng.Find.Execute('JILES, D. C. Introduction to Magnetism and Magnetic Materials, CRC Press, 1998.' + $nl + 'COEY, J. M. D. Magnetism and Magnetic Materials, Cambridge University Press, 2010.' + $nl + 'BUSCHOW, K. H. J.; DE BOER, F. R. Physics of Magnetism and Magnetic Materials, Springer, 2003.' + $nl + 'CULLITY, B. D.; GRAHAM, C. D. Introduction to Magnetic Materials, Wiley-IEEE Press, 2008.' + $nl + 'POOLE, C. P. et al., Superconductivity, Academic Press, 2007.' + $nl + 'SHEAHEN, T. P. Introduction to High-Temperature Superconductivity, Kluwer Academic, 2002.' + $nl + 'LEE, P. J. Engineering Superconductivity, Wiley-IEEE Press, 2001.', $true, $false, $false, $false, $false, $true, 1, $false, '5840726 - Cristina Bormio Nunes', 2)
if (-not $ok) { Write-Output "FAILED: para 19 run 0" } else { Write-Output "OK: para 19 run 0" }
